$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4285.7144
$ws.Cells.Item(17, 10).Value = 5249.75
$ws.Cells.Item(17, 12).Value = 15749.25
$ws.Cells.Item(17, 14).Value = -16085.25
$ws.Cells.Item(33, 8).Value = 1193.3846
$ws.Cells.Item(33, 9).Value = 390.1111
$ws.Cells.Item(33, 11).Value = 390.1111
$ws.Cells.Item(33, 13).Value = -161.1111
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).Value = ""
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).Value = ""
$ws.Cells.Item(88, 8).Value = 1512.5
$ws.Cells.Item(88, 10).Value = 2252
$ws.Cells.Item(88, 12).Value = 2252
$ws.Cells.Item(88, 14).Value = -3064
$ws.Cells.Item(91, 8).Value = 1512.5
$ws.Cells.Item(91, 10).Value = 2252
$ws.Cells.Item(91, 12).Value = 2252
$ws.Cells.Item(91, 14).Value = -5060
$ws.Cells.Item(98, 8).Value = 3576.0715
$ws.Cells.Item(98, 9).Value = 3468
$ws.Cells.Item(98, 10).Value = 3846.25
$ws.Cells.Item(98, 11).Value = 3468
$ws.Cells.Item(98, 12).Value = 3846.25
$ws.Cells.Item(98, 13).Value = -1970
$ws.Cells.Item(98, 14).Value = -6842.25
$ws.Cells.Item(122, 8).Value = 3576.0715
$ws.Cells.Item(122, 9).Value = 3468
$ws.Cells.Item(122, 10).Value = 3846.25
$ws.Cells.Item(122, 11).Value = 10404
$ws.Cells.Item(122, 12).Value = 11538.75
$ws.Cells.Item(122, 13).Value = -7954
$ws.Cells.Item(122, 14).Value = -16438.75
$ws.Cells.Item(132, 8).Value = 2659.4375
$ws.Cells.Item(132, 9).Value = 2778.4546
$ws.Cells.Item(132, 10).Value = 2397.6
$ws.Cells.Item(132, 11).Value = 8335.363799999999
$ws.Cells.Item(132, 12).Value = 7192.799999999999
$ws.Cells.Item(132, 13).Value = -5805.363799999999
$ws.Cells.Item(132, 14).Value = -12252.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 3799.6667
$ws.Cells.Item(12, 9).Value = 699.5
$ws.Cells.Item(12, 11).Value = 699.5
$ws.Cells.Item(12, 13).Value = -526.5
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 14).Value = ""
$ws.Cells.Item(61, 8).Value = 3714.8
$ws.Cells.Item(61, 9).Value = 2999.6667
$ws.Cells.Item(61, 11).Value = 2999.6667
$ws.Cells.Item(61, 13).Value = -2787.6667
$ws.Cells.Item(74, 8).Value = 3519.3333
$ws.Cells.Item(74, 9).Value = 3584.375
$ws.Cells.Item(74, 11).Value = 3584.375
$ws.Cells.Item(74, 13).Value = -2710.375
$ws.Cells.Item(77, 8).Value = 3519.3333
$ws.Cells.Item(77, 9).Value = 3584.375
$ws.Cells.Item(77, 11).Value = 17921.875
$ws.Cells.Item(77, 13).Value = -13553.875
$ws.Cells.Item(88, 8).Value = 4700.7144
$ws.Cells.Item(88, 10).Value = 5645.6
$ws.Cells.Item(88, 12).Value = 5645.6
$ws.Cells.Item(88, 14).Value = -6457.6
$ws.Cells.Item(91, 8).Value = 4700.7144
$ws.Cells.Item(91, 10).Value = 5645.6
$ws.Cells.Item(91, 12).Value = 5645.6
$ws.Cells.Item(91, 14).Value = -8453.6
$ws.Cells.Item(122, 8).Value = 3486.4285
$ws.Cells.Item(122, 9).Value = 2734.1667
$ws.Cells.Item(122, 11).Value = 8202.500100000001
$ws.Cells.Item(122, 13).Value = -5752.500100000001
$ws.Cells.Item(132, 8).Value = 2889
$ws.Cells.Item(132, 9).Value = 2965.5557
$ws.Cells.Item(132, 10).Value = 2200
$ws.Cells.Item(132, 11).Value = 8896.667099999999
$ws.Cells.Item(132, 12).Value = 6600
$ws.Cells.Item(132, 13).Value = -6366.667099999999
$ws.Cells.Item(132, 14).Value = -11660
$ws.Cells.Item(136, 8).Value = 3714.8
$ws.Cells.Item(136, 9).Value = 2999.6667
$ws.Cells.Item(136, 11).Value = 8999.000100000001
$ws.Cells.Item(136, 13).Value = -6449.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(30, 8).Value = 1500
$ws.Cells.Item(30, 10).Value = 1500
$ws.Cells.Item(30, 12).Value = 1500
$ws.Cells.Item(30, 14).Value = -1750
$ws.Cells.Item(86, 8).Value = 3199
$ws.Cells.Item(86, 9).Value = 2982.3333
$ws.Cells.Item(86, 10).Value = 4499
$ws.Cells.Item(86, 11).Value = 2982.3333
$ws.Cells.Item(86, 12).Value = 4499
$ws.Cells.Item(86, 13).Value = -1859.3333
$ws.Cells.Item(86, 14).Value = -6745
$ws.Cells.Item(89, 8).Value = 3199
$ws.Cells.Item(89, 9).Value = 2982.3333
$ws.Cells.Item(89, 10).Value = 4499
$ws.Cells.Item(89, 11).Value = 14911.6665
$ws.Cells.Item(89, 12).Value = 22495
$ws.Cells.Item(89, 13).Value = -9295.666499999999
$ws.Cells.Item(89, 14).Value = -33727
$ws.Cells.Item(107, 8).Value = 1633
$ws.Cells.Item(107, 10).Value = 1920
$ws.Cells.Item(107, 12).Value = 1920
$ws.Cells.Item(107, 14).Value = -5760

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 372.5
$ws.Cells.Item(22, 9).Value = 350
$ws.Cells.Item(22, 11).Value = 350
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(99, 8).Value = 5068.8335
$ws.Cells.Item(99, 9).Value = 4875
$ws.Cells.Item(99, 11).Value = 4875
$ws.Cells.Item(99, 13).Value = -3377
$ws.Cells.Item(122, 8).Value = 2500
$ws.Cells.Item(122, 9).Value = 3000
$ws.Cells.Item(122, 10).Value = 2250
$ws.Cells.Item(122, 11).Value = 9000
$ws.Cells.Item(122, 12).Value = 6750
$ws.Cells.Item(122, 13).Value = -6550
$ws.Cells.Item(122, 14).Value = -11650
$ws.Cells.Item(126, 8).Value = 5068.8335
$ws.Cells.Item(126, 9).Value = 4875
$ws.Cells.Item(126, 11).Value = 14625
$ws.Cells.Item(126, 13).Value = -12155
$ws.Cells.Item(134, 8).Value = 4002.3635
$ws.Cells.Item(134, 9).Value = 2603.1428
$ws.Cells.Item(134, 10).Value = 6451
$ws.Cells.Item(134, 11).Value = 7809.428400000001
$ws.Cells.Item(134, 12).Value = 19353
$ws.Cells.Item(134, 13).Value = -5274.428400000001
$ws.Cells.Item(134, 14).Value = -24423

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 4465.143
$ws.Cells.Item(6, 9).Value = 209.33333
$ws.Cells.Item(6, 11).Value = 627.99999
$ws.Cells.Item(6, 13).Value = -514.99999
$ws.Cells.Item(137, 8).Value = 3120
$ws.Cells.Item(137, 9).Value = 950
$ws.Cells.Item(137, 11).Value = 2850
$ws.Cells.Item(137, 13).Value = 2250

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 3717.8333
$ws.Cells.Item(9, 9).Value = 461.4
$ws.Cells.Item(9, 10).Value = 20000
$ws.Cells.Item(9, 11).Value = 461.4
$ws.Cells.Item(9, 12).Value = 20000
$ws.Cells.Item(9, 13).Value = -291.4
$ws.Cells.Item(9, 14).Value = -20340
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).Value = ""
$ws.Cells.Item(70, 8).Value = 3994.3333
$ws.Cells.Item(70, 9).Value = 3994.3333
$ws.Cells.Item(70, 11).Value = 3994.3333
$ws.Cells.Item(70, 13).Value = -3724.3333
$ws.Cells.Item(73, 8).Value = 3994.3333
$ws.Cells.Item(73, 9).Value = 3994.3333
$ws.Cells.Item(73, 11).Value = 3994.3333
$ws.Cells.Item(73, 13).Value = -3058.3333
$ws.Cells.Item(80, 8).Value = 3785.7144
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 3785.7144
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 3785.7144
$ws.Cells.Item(80, 13).Value = ""
$ws.Cells.Item(80, 14).Value = -5781.7144
$ws.Cells.Item(83, 8).Value = 3785.7144
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 3785.7144
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 18928.572
$ws.Cells.Item(83, 13).Value = ""
$ws.Cells.Item(83, 14).Value = -28912.572
$ws.Cells.Item(102, 8).Value = 1196.9412
$ws.Cells.Item(102, 9).Value = 934.46155
$ws.Cells.Item(102, 10).Value = 2050
$ws.Cells.Item(102, 11).Value = 934.46155
$ws.Cells.Item(102, 12).Value = 2050
$ws.Cells.Item(102, 13).Value = 687.53845
$ws.Cells.Item(102, 14).Value = -5294
$ws.Cells.Item(122, 8).Value = 13907979
$ws.Cells.Item(122, 10).Value = 23481.6
$ws.Cells.Item(122, 12).Value = 70444.79999999999
$ws.Cells.Item(122, 14).Value = -75344.79999999999
$ws.Cells.Item(126, 8).Value = 3858
$ws.Cells.Item(126, 9).Value = 1003
$ws.Cells.Item(126, 11).Value = 3009
$ws.Cells.Item(126, 13).Value = -539

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1100
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).Value = ""
$ws.Cells.Item(46, 8).Value = 3129.85
$ws.Cells.Item(46, 9).Value = 1966.5555
$ws.Cells.Item(46, 11).Value = 1966.5555
$ws.Cells.Item(46, 13).Value = -1778.5555
$ws.Cells.Item(126, 8).Value = 1100
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).Value = ""
$ws.Cells.Item(132, 8).Value = 9506.157999999999
$ws.Cells.Item(132, 9).Value = 9726.125
$ws.Cells.Item(132, 11).Value = 29178.375
$ws.Cells.Item(132, 13).Value = -26648.375
$ws.Cells.Item(136, 8).Value = 2819.0667
$ws.Cells.Item(136, 9).Value = 2698.6428
$ws.Cells.Item(136, 11).Value = 8095.928400000001
$ws.Cells.Item(136, 13).Value = -5545.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 1002
$ws.Cells.Item(8, 9).Value = 1000
$ws.Cells.Item(8, 10).Value = 1004
$ws.Cells.Item(8, 11).Value = 1000
$ws.Cells.Item(8, 12).Value = 1004
$ws.Cells.Item(8, 13).Value = -860
$ws.Cells.Item(8, 14).Value = -1284
$ws.Cells.Item(64, 8).Value = 90000
$ws.Cells.Item(64, 9).Value = 90000
$ws.Cells.Item(64, 11).Value = 90000
$ws.Cells.Item(64, 13).Value = -89752
$ws.Cells.Item(67, 8).Value = 90000
$ws.Cells.Item(67, 9).Value = 90000
$ws.Cells.Item(67, 11).Value = 90000
$ws.Cells.Item(67, 13).Value = -89142
$ws.Cells.Item(132, 8).Value = 1084.7858
$ws.Cells.Item(132, 9).Value = 968.2
$ws.Cells.Item(132, 11).Value = 2904.6
$ws.Cells.Item(132, 13).Value = -374.6000000000004
$ws.Cells.Item(136, 8).Value = 1837.5264
$ws.Cells.Item(136, 9).Value = 1717.3889
$ws.Cells.Item(136, 11).Value = 5152.1667
$ws.Cells.Item(136, 13).Value = -2602.1667
